$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three existing hyperlinks (B3, B4, B5) get rebuilt below with the new
# row layout - the B5 link moves up to B2, so clear them all first.
$ws.Range("A1:F5").Hyperlinks.Delete()

# Reset A2:D5 to the plain "Normal" style first so stale per-cell formatting
# (B5 previously carried the blue/underlined Hyperlink font) doesn't leak
# into the fresh layout below.
$ws.Range("A2:D5").Style = "Normal"

# --- Formatting first: thin box border around A2:D5, column B as text.
#     (Number format must be applied before the values are written so the
#     numeric-looking strings "0816555895"/"86772272121" keep their literal
#     text form - incl. the leading zero - instead of becoming doubles.) ---
$ws.Range("A2:D5").Borders.LineStyle = 1
$ws.Range("A2:D5").Borders.Weight = 2
$ws.Range("B2:B5").NumberFormat = "@"

# --- Row 2: falih hilmi / 0816555895 / admin123 / Login Success Go To page Profile ---
$ws.Range("A2").Value = "falih hilmi"
$ws.Range("C2").Value = "admin123"
$ws.Range("D2").Value = "Login Success Go To page Profile"

# --- Row 3: Rizki Maulana / !@#$%^&**&^%# / TestingQA123 / Akun Anda ... berbelanja. ---
$ws.Range("A3").Value = "Rizki Maulana"
$ws.Range("B3").Value = '!@#$%^&**&^%#'
$ws.Range("C3").Value = "TestingQA123"
$ws.Range("D3").Value = "Akun Anda untuk sementara tidak dapat digunakan untuk berbelanja."

# --- Row 4: falih hilmi / hilmi.falih@yahoo.com / admin123 / Login Success Go To page Profile ---
$ws.Range("A4").Value = "falih hilmi"
$ws.Range("B4").Value = "hilmi.falih@yahoo.com"
$ws.Range("C4").Value = "admin123"
$ws.Range("D4").Value = "Login Success Go To page Profile"

# --- Row 5: Cynthia / 86772272121 / 123admin / Invalid email/phone number or password (no hyperlink anymore) ---
$ws.Range("A5").Value = "Cynthia"
$ws.Range("B5").Value = "86772272121"
$ws.Range("C5").Value = "123admin"
$ws.Range("D5").Value = "Invalid email/phone number or password"

# --- Hyperlinks rebuilt in order so relationship ids line up: rId1->B3, rId2->B4, rId3->B2 ---
$ws.Hyperlinks.Add($ws.Range("B3"), 'mailto:!@#$%^&**&^%#')
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:hilmi.falih@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("B2"), 'mailto:!@#$%^&**&^%#', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "hilmi.falih@yahoo.com")

# Hyperlinks.Add() stamps its TextToDisplay text into the cell itself; put
# B2's real value back afterwards (the hyperlink keeps its stale "display"
# attribute pointing at the old email text, same as in the source workbook).
$ws.Range("B2").Value = "0816555895"

# --- Column C is wider now that it holds "TestingQA123" ---
$ws.Columns.Item(3).AutoFit()

# --- Selection moves to D14 in the saved view state ---
$ws.Range("D14").Select()
